# ---------------------------------------------------------------------------
# "Update countries & provincias Spain"
#
# 1) sharedStrings reorder (country list nav order):
#      - "Mozambique" moves to just before "Republica de Yibuti" (rows 112/113)
#      - "Jordania" moves to just before "Tailandia" (rows 128-134 ripple down)
# 2) Refreshed case counts (new scrape) for a batch of countries
# 3) Updated "last refreshed" timestamp banner
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Banner: "Datos actualizados a ..." timestamp (A1) -----------------------
$ws.Range("A1").Value = "Datos actualizados a 14 de Septiembre de 2020 a las 20:14"

# --- Country re-ordering: Mozambique now sorts before Republica de Yibuti ---
# --- and Jordania now sorts before Tailandia; rows get rewritten in place ---
# --- with the refreshed totals that came with the new ranking. --------------

# Row 112: Mozambique
$row = New-Object 'object[,]' 1,8
$row[0,0] = "Mozambique"
$row[0,1] = 5482
$row[0,2] = 213
$row[0,3] = 3024
$row[0,4] = 2423
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 35
$ws.Range("A112:H112").Value = $row

# Row 113: Republica de Yibuti
$row = New-Object 'object[,]' 1,8
$row[0,0] = "Republica de Yibuti"
$row[0,1] = 5396
$row[0,2] = 1
$row[0,3] = 5331
$row[0,4] = 4
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 61
$ws.Range("A113:H113").Value = $row

# Row 128: Jordania
$row = New-Object 'object[,]' 1,8
$row[0,0] = "Jordania"
$row[0,1] = 3528
$row[0,2] = 214
$row[0,3] = 2255
$row[0,4] = 1247
$row[0,5] = 0
$row[0,6] = 2
$row[0,7] = 26
$ws.Range("A128:H128").Value = $row

# Row 129: Tailandia
$row = New-Object 'object[,]' 1,8
$row[0,0] = "Tailandia"
$row[0,1] = 3475
$row[0,2] = 2
$row[0,3] = 3312
$row[0,4] = 105
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 58
$ws.Range("A129:H129").Value = $row

# Row 130: Gambia
$row = New-Object 'object[,]' 1,8
$row[0,0] = "Gambia"
$row[0,1] = 3405
$row[0,2] = 0
$row[0,3] = 1723
$row[0,4] = 1579
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 103
$ws.Range("A130:H130").Value = $row

# Row 131: Somalia
$row = New-Object 'object[,]' 1,8
$row[0,0] = "Somalia"
$row[0,1] = 3389
$row[0,2] = 0
$row[0,3] = 2803
$row[0,4] = 488
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 98
$ws.Range("A131:H131").Value = $row

# Row 132: Angola
$row = New-Object 'object[,]' 1,8
$row[0,0] = "Angola"
$row[0,1] = 3388
$row[0,2] = 0
$row[0,3] = 1301
$row[0,4] = 1953
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 134
$ws.Range("A132:H132").Value = $row

# Row 133: Lituania
$row = New-Object 'object[,]' 1,8
$row[0,0] = "Lituania"
$row[0,1] = 3386
$row[0,2] = 51
$row[0,3] = 2071
$row[0,4] = 1228
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 87
$ws.Range("A133:H133").Value = $row

# Row 134: Mayotte
$row = New-Object 'object[,]' 1,8
$row[0,0] = "Mayotte"
$row[0,1] = 3374
$row[0,2] = 0
$row[0,3] = 2964
$row[0,4] = 370
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 40
$ws.Range("A134:H134").Value = $row

# --- Refreshed counts only (country/rank unchanged) -------------------------

# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 6723298
$ws.Cells.Item(4, 3).Value = 14840
$ws.Cells.Item(4, 5).Value = 2541827
$ws.Cells.Item(4, 7).Value = 176
$ws.Cells.Item(4, 8).Value = 198696

# Row 5: India
$ws.Cells.Item(5, 2).Value = 4926734
$ws.Cells.Item(5, 3).Value = 81731
$ws.Cells.Item(5, 4).Value = 3855983
$ws.Cells.Item(5, 5).Value = 989924
$ws.Cells.Item(5, 7).Value = 1073
$ws.Cells.Item(5, 8).Value = 80827

# Row 16: Francia
$ws.Cells.Item(16, 2).Value = 387252
$ws.Cells.Item(16, 3).Value = 6158
$ws.Cells.Item(16, 4).Value = 89507
$ws.Cells.Item(16, 5).Value = 266795
$ws.Cells.Item(16, 7).Value = 34
$ws.Cells.Item(16, 8).Value = 30950

# Row 27: Israel
$ws.Cells.Item(27, 2).Value = 159290
$ws.Cells.Item(27, 3).Value = 3686
$ws.Cells.Item(27, 4).Value = 118547
$ws.Cells.Item(27, 5).Value = 39607
$ws.Cells.Item(27, 7).Value = 17
$ws.Cells.Item(27, 8).Value = 1136

# Row 41: Marruecos
$ws.Cells.Item(41, 2).Value = 88203
$ws.Cells.Item(41, 3).Value = 1517
$ws.Cells.Item(41, 4).Value = 68970
$ws.Cells.Item(41, 5).Value = 17619
$ws.Cells.Item(41, 7).Value = 36
$ws.Cells.Item(41, 8).Value = 1614

# Row 73: Irlanda
$ws.Cells.Item(73, 2).Value = 31192
$ws.Cells.Item(73, 3).Value = 207
$ws.Cells.Item(73, 5).Value = 6044

# Row 77: Libano
$ws.Cells.Item(77, 2).Value = 24857
$ws.Cells.Item(77, 3).Value = 547
$ws.Cells.Item(77, 5).Value = 15846

# Row 94: Albania
$ws.Cells.Item(94, 2).Value = 11520
$ws.Cells.Item(94, 3).Value = 167
$ws.Cells.Item(94, 4).Value = 6615
$ws.Cells.Item(94, 5).Value = 4567
$ws.Cells.Item(94, 7).Value = 4
$ws.Cells.Item(94, 8).Value = 338

# Row 135: Sri Lanka
$ws.Cells.Item(135, 2).Value = 3262
$ws.Cells.Item(135, 3).Value = 28
$ws.Cells.Item(135, 5).Value = 244

# Row 161: Principado de Andorra
$ws.Cells.Item(161, 2).Value = 1438
$ws.Cells.Item(161, 3).Value = 94
$ws.Cells.Item(161, 4).Value = 945
$ws.Cells.Item(161, 5).Value = 440
